{"js": "// Replace the 25 division problems in the practice table with their\n// updated values, matching the cells in row-major (top-to-bottom,\n// left-to-right) order.\n\nconst table = context.document.body.tables.getFirst();\n\n// Only every 4th row (0, 4, 8, 12, 16) actually holds problem text; the\n// rows in between are blank spacer rows that keep their own layout.\nconst rowIndexes = [0, 4, 8, 12, 16];\n\nconst replacements = [\n  [\"93\u00f76=\", \"48\u00f75=\"],\n  [\"45\u00f74=\", \"37\u00f78=\"],\n  [\"11\u00f72=\", \"10\u00f74=\"],\n  [\"45\u00f78=\", \"62\u00f72=\"],\n  [\"36\u00f78=\", \"30\u00f79=\"],\n  [\"96\u00f77=\", \"10\u00f73=\"],\n  [\"43\u00f76=\", \"98\u00f77=\"],\n  [\"83\u00f79=\", \"92\u00f78=\"],\n  [\"93\u00f72=\", \"72\u00f76=\"],\n  [\"26\u00f74=\", \"36\u00f79=\"],\n  [\"63\u00f74=\", \"10\u00f77=\"],\n  [\"36\u00f72=\", \"85\u00f73=\"],\n  [\"57\u00f76=\", \"10\u00f72=\"],\n  [\"50\u00f77=\", \"97\u00f79=\"],\n  [\"70\u00f78=\", \"51\u00f77=\"],\n  [\"26\u00f77=\", \"35\u00f72=\"],\n  [\"92\u00f75=\", \"36\u00f74=\"],\n  [\"15\u00f73=\", \"23\u00f74=\"],\n  [\"39\u00f74=\", \"44\u00f72=\"],\n  [\"62\u00f79=\", \"24\u00f73=\"],\n  [\"46\u00f78=\", \"46\u00f76=\"],\n  [\"91\u00f75=\", \"76\u00f75=\"],\n  [\"33\u00f74=\", \"92\u00f75=\"],\n  [\"11\u00f74=\", \"33\u00f79=\"],\n  [\"27\u00f72=\", \"26\u00f74=\"],\n];\n\nconst cells = [];\nfor (const rowIndex of rowIndexes) {\n  for (let col = 0; col < 5; col++) {\n    cells.push(table.getCell(rowIndex, col));\n  }\n}\n\nfor (const cell of cells) {\n  cell.load(\"value\");\n}\nawait context.sync();\n\nfor (let i = 0; i < cells.length; i++) {\n  const [expected, updated] = replacements[i];\n  const currentText = cells[i].value.trim();\n  if (currentText !== expected) {\n    throw new Error(\n      `Unexpected cell text at index ${i}: got \"${currentText}\", expected \"${expected}\"`\n    );\n  }\n  cells[i].value = updated;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division problems in the practice table with their\n# updated values, matching the cells in row-major (top-to-bottom,\n# left-to-right) order. COM collections are 1-based.\n\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n# Only every 4th row (1, 5, 9, 13, 17) actually holds problem text; the\n# rows in between are blank spacer rows that keep their own layout.\n$rowIndexes = @(1, 5, 9, 13, 17)\n\n$replacements = @(\n  @(\"93\u00f76=\", \"48\u00f75=\"),\n  @(\"45\u00f74=\", \"37\u00f78=\"),\n  @(\"11\u00f72=\", \"10\u00f74=\"),\n  @(\"45\u00f78=\", \"62\u00f72=\"),\n  @(\"36\u00f78=\", \"30\u00f79=\"),\n  @(\"96\u00f77=\", \"10\u00f73=\"),\n  @(\"43\u00f76=\", \"98\u00f77=\"),\n  @(\"83\u00f79=\", \"92\u00f78=\"),\n  @(\"93\u00f72=\", \"72\u00f76=\"),\n  @(\"26\u00f74=\", \"36\u00f79=\"),\n  @(\"63\u00f74=\", \"10\u00f77=\"),\n  @(\"36\u00f72=\", \"85\u00f73=\"),\n  @(\"57\u00f76=\", \"10\u00f72=\"),\n  @(\"50\u00f77=\", \"97\u00f79=\"),\n  @(\"70\u00f78=\", \"51\u00f77=\"),\n  @(\"26\u00f77=\", \"35\u00f72=\"),\n  @(\"92\u00f75=\", \"36\u00f74=\"),\n  @(\"15\u00f73=\", \"23\u00f74=\"),\n  @(\"39\u00f74=\", \"44\u00f72=\"),\n  @(\"62\u00f79=\", \"24\u00f73=\"),\n  @(\"46\u00f78=\", \"46\u00f76=\"),\n  @(\"91\u00f75=\", \"76\u00f75=\"),\n  @(\"33\u00f74=\", \"92\u00f75=\"),\n  @(\"11\u00f74=\", \"33\u00f79=\"),\n  @(\"27\u00f72=\", \"26\u00f74=\")\n)\n\n$idx = 0\nforeach ($r in $rowIndexes) {\n  for ($c = 1; $c -le 5; $c++) {\n    $cell = $t.Cell($r, $c)\n    $range = $cell.Range\n    $expected = $replacements[$idx][0]\n    $updated = $replacements[$idx][1]\n\n    $currentText = $range.Text.Substring(0, $expected.Length)\n    if ($currentText -ne $expected) {\n      throw \"Unexpected cell text at row $r, col $c`: got '$currentText', expected '$expected'\"\n    }\n\n    $range.Text = $updated\n    $idx++\n  }\n}\n"}
